# Generate Report for Handoff
#
# The localization-status report is regenerated: the row for
# 1bbcfdae-0611-42d0-b919-260575de5ba1.md now sorts before the row for
# 189213fe-347a-4762-9be1-2ae0b57be902.md (rows 2 and 3 swap places across
# all three sheets), and 189213fe-...md has moved from "In Translation" to
# "Ready for handoff" with a fresh handoff file/timestamp, while
# 1bbcfdae-...md stays "In Translation".

$wb = $excel.ActiveWorkbook

function Set-CellAndHyperlink($ws, [string]$addr, [string]$newValue, [bool]$hasHyperlink) {
    $ws.Range($addr).Value = $newValue
    if ($hasHyperlink) {
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range.Address(0, 0) -eq $addr) {
                $h.TextToDisplay = $newValue
            }
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

Set-CellAndHyperlink $ws "A2" "1bbcfdae-0611-42d0-b919-260575de5ba1.md" $true
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"

Set-CellAndHyperlink $ws "A3" "189213fe-347a-4762-9be1-2ae0b57be902.md" $true
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

Set-CellAndHyperlink $ws "A2" "1bbcfdae-0611-42d0-b919-260575de5ba1.md" $true
$ws.Range("B2").Value = "In Translation"
Set-CellAndHyperlink $ws "C2" "1bbcfdae-0611-42d0-b919-260575de5ba1.fcf8607e21a11b0e90c0aa0e3e35d21678cc121a.zh-cn.xlf" $true
$ws.Range("D2").Value = "2016-03-11 00:16:26"

Set-CellAndHyperlink $ws "A3" "189213fe-347a-4762-9be1-2ae0b57be902.md" $true
$ws.Range("B3").Value = "Ready for handoff"
Set-CellAndHyperlink $ws "C3" "189213fe-347a-4762-9be1-2ae0b57be902.f4c13c59ebad7e5e1c38915a39fd2517f81dfdd8.zh-cn.xlf" $true
$ws.Range("D3").Value = "2016-03-11 00:20:19"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

Set-CellAndHyperlink $ws "A2" "1bbcfdae-0611-42d0-b919-260575de5ba1.md" $true
$ws.Range("B2").Value = "In Translation"
Set-CellAndHyperlink $ws "C2" "1bbcfdae-0611-42d0-b919-260575de5ba1.fcf8607e21a11b0e90c0aa0e3e35d21678cc121a.de-de.xlf" $true
$ws.Range("D2").Value = "2016-03-11 00:17:43"

Set-CellAndHyperlink $ws "A3" "189213fe-347a-4762-9be1-2ae0b57be902.md" $true
$ws.Range("B3").Value = "Ready for handoff"
Set-CellAndHyperlink $ws "C3" "189213fe-347a-4762-9be1-2ae0b57be902.f4c13c59ebad7e5e1c38915a39fd2517f81dfdd8.de-de.xlf" $true
$ws.Range("D3").Value = "2016-03-11 00:20:25"
